$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1450048780487805
$ws.Range("V2").Value = 0.0002448603057459146
$ws.Range("Z2").Value = -0.2462948787810065
$ws.Range("AB2").Value = -1005.85874068368
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -1005.85874068368

# Row 3
$ws.Range("T3").Value = 0.1492487804878049
$ws.Range("V3").Value = 0.0001488973818309612
$ws.Range("Z3").Value = -0.2680607270169032
$ws.Range("AB3").Value = -1800.305174749309
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -1800.305174749309

# Row 4
$ws.Range("T4").Value = 0.1469268292682927
$ws.Range("V4").Value = 0.0002222807942365138
$ws.Range("Z4").Value = -0.1747274048841794
$ws.Range("AB4").Value = -786.0661353326996
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -786.0661353326996

# Row 5
$ws.Range("T5").Value = 0.1418926829268293
$ws.Range("V5").Value = 0.0002529432437181515
$ws.Range("Z5").Value = -0.2668150213665862
$ws.Range("AB5").Value = -1054.841463422884
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -1054.841463422884

# Row 6
$ws.Range("T6").Value = 0.1446439024390244
$ws.Range("V6").Value = 0.0001851607801792304
$ws.Range("Z6").Value = -0.3043911547471467
$ws.Range("AB6").Value = -1643.928884143309
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -1643.928884143309

# Row 7
$ws.Range("T7").Value = 0.1429268292682927
$ws.Range("V7").Value = 0.0003232296608680373
$ws.Range("Z7").Value = -0.2232340414394441
$ws.Range("AB7").Value = -690.6360042576111
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -690.6360042576111

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.0004493416693347922
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"

# Row 9
$ws.Range("T9").Value = 0.1450048780487805
$ws.Range("V9").Value = 0.0002448603057459146
$ws.Range("Z9").Value = 0.3051536134507343
$ws.Range("AB9").Value = 1246.235532219683
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 1246.235532219683

# Row 10
$ws.Range("T10").Value = 0.1492487804878049
$ws.Range("V10").Value = 0.0001488973818309612
$ws.Range("Z10").Value = 0.3333869806534889
$ws.Range("AB10").Value = 2239.038568401245
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 2239.038568401245

# Row 11
$ws.Range("T11").Value = 0.1469268292682927
$ws.Range("V11").Value = 0.0002222807942365138
$ws.Range("Z11").Value = 0.18021753262247
$ws.Range("AB11").Value = 810.765200122116
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = 810.765200122116

# Row 12
$ws.Range("T12").Value = 0.1418926829268293
$ws.Range("V12").Value = 0.0002529432437181515
$ws.Range("Z12").Value = 0.3197094174206862
$ws.Range("AB12").Value = 1263.9571341029
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 1263.9571341029

# Row 13
$ws.Range("T13").Value = 0.1446439024390244
$ws.Range("V13").Value = 0.0001851607801792304
$ws.Range("Z13").Value = 0.3933418426589942
$ws.Range("AB13").Value = 2124.325908965443
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 2124.325908965443

# Row 14
$ws.Range("T14").Value = 0.1429268292682927
$ws.Range("V14").Value = 0.0003232296608680373
$ws.Range("Z14").Value = 0.2726768320113346
$ws.Range("AB14").Value = 843.6008975137293
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 843.6008975137293

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = 0.0001641955856056655
$ws.Range("AB15").Value = "Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "Inf"

